$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Comput_time (column J) values
$ws.Range("J2").Value = 104.198
$ws.Range("J3").Value = 6.22
$ws.Range("J4").Value = 26.242
$ws.Range("J5").Value = 3.072
$ws.Range("J6").Value = 210.918
$ws.Range("J7").Value = 11.879
$ws.Range("J8").Value = 9.888999999999999
$ws.Range("J9").Value = 15.589
$ws.Range("J10").Value = 4.776
$ws.Range("J11").Value = 2.244
$ws.Range("J12").Value = 9.765000000000001
$ws.Range("J13").Value = 15.772
$ws.Range("J14").Value = 6.143

# Row 15 and 16: Test_color changes from "Orange (not order 2)" to "Green", plus Comput_time
$ws.Range("I15").Value = "Green"
$ws.Range("J15").Value = 62.615

$ws.Range("I16").Value = "Green"
$ws.Range("J16").Value = 3.962
